# Mise à jour de l'application
# Adds a new day column (DK) of attendance data to the tracking sheet,
# mirroring the style of the existing last day column (DJ), and updates
# the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date for the added day: 2026-02-04 (serial 46057), one day after
# the previous last column (DJ1 = 46056 = 2026-02-03).
$newDateCol = 115   # column DK
$lastDateCol = 114  # column DJ

# --- Row 1 (header / date row): set the date value first, then copy DJ1's
# format onto DK1 (format-only paste so the cached formula recalc below still
# sees the new value as a fresh edit).
$ws.Cells.Item(1, $newDateCol).Value = 46057
$ws.Cells.Item(1, $lastDateCol).Copy() | Out-Null
$ws.Cells.Item(1, $newDateCol).PasteSpecial(-4122) | Out-Null

# Attendance value for the new day, per player row. "P" = present, "B" = away
# ("Blessure"), $null = leave the new cell blank (still formatted), and rows
# absent from this map are not touched at all (their data doesn't extend
# that far right, same as before the edit).
$attendance = @{
    2  = "P";  3  = "P";  4  = "P";  5  = "P";  6  = "B";  7  = "P";
    8  = "P";  9  = "P"; 10  = "P"; 11  = "P"; 13  = "P"; 14  = "P";
    15 = "P"; 16  = $null; 17 = $null; 18 = "P"; 19 = "P"; 20 = "P";
    22 = "P"; 24 = "P"; 25 = $null; 26 = "P"; 27 = "P"; 28 = "P";
    29 = "P"; 30 = "P"; 31 = "P"
}

foreach ($row in 2..31) {
    if (-not $attendance.ContainsKey($row)) {
        continue
    }

    # Set the value (if any) BEFORE copying formatting over: a format-only
    # paste performed after the value write can otherwise leave the
    # recalculation engine thinking the cell/its dependents are unchanged.
    $value = $attendance[$row]
    if ($null -ne $value) {
        $ws.Cells.Item($row, $newDateCol).Value = $value
    }

    # Copy the formatting of the row's current last day cell (DJ) onto the
    # new day cell (DK) so the new column matches the existing style exactly.
    $ws.Cells.Item($row, $lastDateCol).Copy() | Out-Null
    $ws.Cells.Item($row, $newDateCol).PasteSpecial(-4122) | Out-Null
}

# Recalculate so every COUNTA/COUNTIF summary column picks up the new day.
$excel.Calculate() | Out-Null

# Move/refresh the active selection like the author's last interaction.
$ws.Range("DN19").Select() | Out-Null
